# Commit: "removed mentions of coverage"
#
# This edit removes the "targetable" and "coverage" columns (and their
# header cells) from the "Model parameters" sheet. This also removes the
# now-unused "targetable" shared string, which causes every other shared
# string that used to live after it to shift its index down by one
# (handled automatically by the engine on save).
#
# It also changes which sheet/cell is active/selected: the previously
# active sheet ("Data constants") becomes inactive, and "Model parameters"
# (with cell I4 selected) becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Model parameters")

# Make "Model parameters" the active sheet and select what was cell I4
# (in the pre-edit column layout) *before* removing the columns, so the
# stored selection reference ends up as I4 (matching the edit author's
# recorded selection, made before they deleted the columns).
$ws.Activate()
$ws.Range("I4").Select()

# Column G holds "targetable" and column H holds "coverage" (by header).
# Delete both whole columns; what was column I ("fromdata") becomes the
# new column G, and what was column J (blank) becomes the new column H.
$ws.Columns.Item(7).EntireColumn.Delete()
$ws.Columns.Item(7).EntireColumn.Delete()
